$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R2").Value = 1.52777777816358
$ws.Range("AI2").Value = 1.52777777816358
$ws.Range("AZ2").Value = 1.52777777816358
$ws.Range("J4").Value = 2
$ws.Range("AA4").Value = 2
$ws.Range("AR4").Value = 2
$ws.Range("R5").Value = 2.722222218364197
$ws.Range("AI5").Value = 2.722222218364197
$ws.Range("AZ5").Value = 2.722222218364197
$ws.Range("J7").Value = 1.090045753188997
$ws.Range("AA7").Value = 1.090045753188997
$ws.Range("AR7").Value = 1.090045753188997
$ws.Range("J8").Value = 1.090045753188997
$ws.Range("AA8").Value = 1.090045753188997
$ws.Range("AR8").Value = 1.090045753188997
$ws.Range("R11").Value = 1.874125878478165
$ws.Range("AI11").Value = 1.874125878478165
$ws.Range("AZ11").Value = 1.874125878478165
$ws.Range("R19").Value = 2.076388883005401
$ws.Range("AI19").Value = 2.076388883005401
$ws.Range("AZ19").Value = 2.076388883005401
$ws.Range("R20").Value = 1.270833316261574
$ws.Range("AI20").Value = 1.270833316261574
$ws.Range("AZ20").Value = 1.270833316261574
$ws.Range("R33").Value = 1.750000003472222
$ws.Range("AI33").Value = 1.750000003472222
$ws.Range("AZ33").Value = 1.750000003472222
$ws.Range("R34").Value = 2.236111107445987
$ws.Range("AI34").Value = 2.236111107445987
$ws.Range("AZ34").Value = 2.236111107445987
$ws.Range("R37").Value = 2.055555549382716
$ws.Range("AI37").Value = 2.055555549382716
$ws.Range("AZ37").Value = 2.055555549382716
$ws.Range("R38").Value = 2.118055550250772
$ws.Range("AI38").Value = 2.118055550250772
$ws.Range("AZ38").Value = 2.118055550250772
$ws.Range("I44").Value = 2.983213836153504
$ws.Range("Z44").Value = 2.983213836153504
$ws.Range("AQ44").Value = 2.983213836153504
$ws.Range("R45").Value = 2.222222211419753
$ws.Range("AI45").Value = 2.222222211419753
$ws.Range("AZ45").Value = 2.222222211419753
$ws.Range("R46").Value = 2.333333317129629
$ws.Range("AI46").Value = 2.333333317129629
$ws.Range("AZ46").Value = 2.333333317129629
$ws.Range("R50").Value = 1.541666636284722
$ws.Range("AI50").Value = 1.541666636284722
$ws.Range("AZ50").Value = 1.541666636284722
$ws.Range("R52").Value = 1.979166659143518
$ws.Range("AI52").Value = 1.979166659143518
$ws.Range("AZ52").Value = 1.979166659143518
$ws.Range("R61").Value = 2.027777771219136
$ws.Range("AI61").Value = 2.027777771219136
$ws.Range("AZ61").Value = 2.027777771219136
$ws.Range("R62").Value = 2.236111100501543
$ws.Range("AI62").Value = 2.236111100501543
$ws.Range("AZ62").Value = 2.236111100501543
$ws.Range("R64").Value = 1.423611099103009
$ws.Range("AI64").Value = 1.423611099103009
$ws.Range("AZ64").Value = 1.423611099103009
$ws.Range("R74").Value = 2.236111107445987
$ws.Range("AI74").Value = 2.236111107445987
$ws.Range("AZ74").Value = 2.236111107445987
$ws.Range("I77").Value = 1
$ws.Range("Z77").Value = 1
$ws.Range("AQ77").Value = 1
$ws.Range("J80").Value = 1.135068629783496
$ws.Range("AA80").Value = 1.135068629783496
$ws.Range("AR80").Value = 1.135068629783496
$ws.Range("R83").Value = 4.069444438464506
$ws.Range("AI83").Value = 4.069444438464506
$ws.Range("AZ83").Value = 4.069444438464506
$ws.Range("J98").Value = 1.136529254589045
$ws.Range("AA98").Value = 1.136529254589045
$ws.Range("AR98").Value = 1.136529254589045
$ws.Range("I101").Value = 1.35594577204859
$ws.Range("Z101").Value = 1.35594577204859
$ws.Range("AQ101").Value = 1.35594577204859
$ws.Range("J107").Value = 2
$ws.Range("AA107").Value = 2
$ws.Range("AR107").Value = 2
$ws.Range("R108").Value = 2.118055550250772
$ws.Range("AI108").Value = 2.118055550250772
$ws.Range("AZ108").Value = 2.118055550250772
$ws.Range("R109").Value = 3.47222222183642
$ws.Range("AI109").Value = 3.47222222183642
$ws.Range("AZ109").Value = 3.47222222183642
$ws.Range("I110").Value = 2
$ws.Range("Z110").Value = 2
$ws.Range("AQ110").Value = 2
$ws.Range("J115").Value = 1.090045753188997
$ws.Range("AA115").Value = 1.090045753188997
$ws.Range("AR115").Value = 1.090045753188997
$ws.Range("I118").Value = 1.826396820127276
$ws.Range("R118").Value = 1.671328673724876
$ws.Range("Z118").Value = 1.826396820127276
$ws.Range("AI118").Value = 1.671328673724876
$ws.Range("AQ118").Value = 1.826396820127276
$ws.Range("AZ118").Value = 1.671328673724876
$ws.Range("R120").Value = 1.671328673724876
$ws.Range("AI120").Value = 1.671328673724876
$ws.Range("AZ120").Value = 1.671328673724876
$ws.Range("R121").Value = 1.700000002857143
$ws.Range("AI121").Value = 1.700000002857143
$ws.Range("AZ121").Value = 1.700000002857143
$ws.Range("R124").Value = 2.076388883005401
$ws.Range("AI124").Value = 2.076388883005401
$ws.Range("AZ124").Value = 2.076388883005401
$ws.Range("R132").Value = 1.700000002857143
$ws.Range("AI132").Value = 1.700000002857143
$ws.Range("AZ132").Value = 1.700000002857143
$ws.Range("R134").Value = 2.118055550250772
$ws.Range("AI134").Value = 2.118055550250772
$ws.Range("AZ134").Value = 2.118055550250772
$ws.Range("R139").Value = 1.700000002857143
$ws.Range("AI139").Value = 1.700000002857143
$ws.Range("AZ139").Value = 1.700000002857143
$ws.Range("R160").Value = 2.236111100501543
$ws.Range("AI160").Value = 2.236111100501543
$ws.Range("AZ160").Value = 2.236111100501543
$ws.Range("R168").Value = 1.671328673724876
$ws.Range("AI168").Value = 1.671328673724876
$ws.Range("AZ168").Value = 1.671328673724876
$ws.Range("J178").Value = 1.090045753188997
$ws.Range("AA178").Value = 1.090045753188997
$ws.Range("AR178").Value = 1.090045753188997
$ws.Range("J179").Value = 1.135068629783496
$ws.Range("AA179").Value = 1.135068629783496
$ws.Range("AR179").Value = 1.135068629783496
$ws.Range("R197").Value = 2.222222211419753
$ws.Range("AI197").Value = 2.222222211419753
$ws.Range("AZ197").Value = 2.222222211419753
$ws.Range("R199").Value = 2.263888878665123
$ws.Range("AI199").Value = 2.263888878665123
$ws.Range("AZ199").Value = 2.263888878665123
